$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SE1")

# correctif probleme insertion matiere dans bd et mise a jour semestre etudiant lors de l'inscription

# 1) "Groupe" (column D): students moved from semester 2 groups to semester 1 groups (2-A/B/C -> 1-A/B/C)
for ($r = 3; $r -le 63; $r++) {
    $groupe = $ws.Cells.Item($r, 4).Text
    if ($groupe -eq "2-A") { $ws.Cells.Item($r, 4).Value = "1-A" }
    elseif ($groupe -eq "2-B") { $ws.Cells.Item($r, 4).Value = "1-B" }
    elseif ($groupe -eq "2-C") { $ws.Cells.Item($r, 4).Value = "1-C" }
}

# 2) Student number (column A): inscription year corrected from 2015 to 2017 (+20000)
for ($r = 3; $r -le 63; $r++) {
    $num = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $num + 20000
}

# 3) "Moyenne de l'etudiant" (column E): recomputed averages per student
$moyennes = @{
    3 = 9
    4 = 9
    5 = 13
    6 = 20
    7 = 7
    8 = 19
    9 = 19
    10 = 12
    11 = 10
    12 = 13
    13 = 16
    14 = 17
    15 = 11
    16 = 20
    17 = 18
    18 = 9
    19 = 11
    20 = 6
    21 = 5
    22 = 20
    23 = 11
    24 = 17
    25 = 16
    26 = 12
    27 = 20
    28 = 6
    29 = 19
    30 = 6
    31 = 7
    32 = 14
    33 = 17
    34 = 8
    35 = 13
    36 = 7
    37 = 11
    38 = 8
    39 = 17
    40 = 6
    41 = 12
    42 = 9
    43 = 6
    44 = 12
    45 = 7
    46 = 14
    47 = 20
    48 = 17
    49 = 18
    50 = 5
    51 = 19
    52 = 12
    53 = 13
    54 = 9
    55 = 19
    56 = 8
    57 = 14
    58 = 19
    59 = 19
    60 = 7
    61 = 5
    62 = 5
    63 = 17
}
foreach ($r in $moyennes.Keys) {
    $ws.Cells.Item($r, 5).Value = $moyennes[$r]
}
